# Update "想去人数" (want-to-go count) values for two events that appear
# in both the "展览" sheet and the "全部类型" (all types) summary sheet.
#   F2: 3210 -> 3215  (南宁·第二届北极光动漫展)
#   F4: 1060 -> 1069  (南宁·2024良牙动漫秋季盛典（秋典）)

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 3215
    $ws.Range("F4").Value = 1069
}
